$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 159, shifting existing rows 159:189 down to 160:190.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new weekly record
# (same market/category context as its neighbours, new date / quality / prices / origin).
$ws.Cells.Item(159, 1).Value = 5
$ws.Cells.Item(159, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(159, 3).Value = "Maule"
$ws.Cells.Item(159, 4).Value = 44505
$ws.Cells.Item(159, 4).NumberFormat = $ws.Cells.Item(160, 4).NumberFormat
$ws.Cells.Item(159, 5).Value = 7
$ws.Cells.Item(159, 6).Value = 100112006
$ws.Cells.Item(159, 7).Value = "Repollo"
$ws.Cells.Item(159, 8).Value = "Crespo record"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 6000
$ws.Cells.Item(159, 11).Value = 600
$ws.Cells.Item(159, 12).Value = 600
$ws.Cells.Item(159, 13).Value = 600
$ws.Cells.Item(159, 14).Value = "`$/unidad"
$ws.Cells.Item(159, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(159, 16).Value = 600
$ws.Cells.Item(159, 17).Value = 1
$ws.Cells.Item(159, 18).Value = "Hortaliza"
